$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "금주 업무" (this week's task) description in the merged cell B12:C12
# Replace the old task text with the new one describing the procmailrc -> python script connection.
$ws.Range("B12").Value = "procmailrc->python script 연결"

# Move the active selection to B13, matching the cursor position left after the edit.
$ws.Activate()
$ws.Range("B13").Select()
